$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 90: Course, Hours, Notes for the Time Log table
$ws.Range("B90").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("C90").Value = 0.75
$ws.Range("D90").Value = "Finish 1 small problem"

# Update selection to D90 to match the author's final cursor position
$ws.Range("D90").Select()
